# Auto-generated edit script: update crypto price/volume columns (D, E) for rows 2-51
# Ensures the cells stay text (matching the original inlineStr cells) by forcing
# a text NumberFormat before assignment -- otherwise Excel auto-coerces numeric-looking
# strings like "0.9980" or "0.000008797" into floating point numbers / scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.285.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9979"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4481"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3779"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07522"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8934"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.01"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.48"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.766"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.64"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.391"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07123"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9987"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008797"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.13"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.280.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.245"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.97"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.047.55"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.990"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.461"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.384"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08851"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7736"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.191"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.584"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.885"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9983"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.110"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01995"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05315"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.450"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5341"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1731"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.854"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.257"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +14.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.801"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5121"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.709"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.26"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9980"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06380"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.92%  "
